# Updates cryptocurrency Price (column D) and Volume/1h (column E) values
# for the rows whose market data changed, per the upstream GitHub Actions
# scrape refresh. Values are written as text (matching the workbook's
# existing inline-string cell type) so strings such as "332.45" or
# "0.92%" are not silently reinterpreted as numbers/percentages.

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "332.45"
Set-TextValue $ws "E2" "0.92%"
Set-TextValue $ws "D3" "39.25"
Set-TextValue $ws "E3" "-2.40%"
Set-TextValue $ws "D4" "5.755"
Set-TextValue $ws "E4" "3.16%"
Set-TextValue $ws "D5" "0.08040"
Set-TextValue $ws "E5" "-0.93%"
Set-TextValue $ws "D6" "4.501"
Set-TextValue $ws "E6" "-1.08%"
Set-TextValue $ws "D7" "8.636"
Set-TextValue $ws "E7" "-0.43%"
Set-TextValue $ws "D8" "1.962"
Set-TextValue $ws "E8" "-2.32%"
Set-TextValue $ws "D10" "0.9233"
Set-TextValue $ws "E10" "-2.73%"
Set-TextValue $ws "D11" "0.1267"
Set-TextValue $ws "E11" "-0.92%"
Set-TextValue $ws "D12" "0.1948"
Set-TextValue $ws "E12" "-1.43%"
Set-TextValue $ws "D13" "8.748"
Set-TextValue $ws "E13" "16.90%"
Set-TextValue $ws "D14" "0.09227"
Set-TextValue $ws "E14" "-0.04%"
Set-TextValue $ws "D15" "0.03552"
Set-TextValue $ws "E15" "-1.34%"
Set-TextValue $ws "E16" "9.62%"
Set-TextValue $ws "D17" "0.001309"
Set-TextValue $ws "E17" "-0.18%"
Set-TextValue $ws "D18" "0.006187"
Set-TextValue $ws "E18" "-1.96%"
Set-TextValue $ws "E19" "0.11%"
Set-TextValue $ws "E20" "-0.93%"
Set-TextValue $ws "E21" "-0.02%"
Set-TextValue $ws "D22" "0.2668"
Set-TextValue $ws "E22" "7.21%"
Set-TextValue $ws "D23" "0.04447"
Set-TextValue $ws "E23" "0.23%"
Set-TextValue $ws "D24" "0.001264"
Set-TextValue $ws "E24" "2.90%"
Set-TextValue $ws "D25" "0.004485"
Set-TextValue $ws "E25" "4.37%"
Set-TextValue $ws "D26" "0.0001207"
Set-TextValue $ws "E26" "0.44%"
Set-TextValue $ws "D39" "0.02522"
Set-TextValue $ws "E39" "0.10%"
Set-TextValue $ws "D40" "0.05471"
Set-TextValue $ws "E40" "4.42%"
Set-TextValue $ws "D41" "0.007454"
Set-TextValue $ws "E41" "-4.38%"
Set-TextValue $ws "D42" "0.009907"
Set-TextValue $ws "E42" "11.64%"
Set-TextValue $ws "D43" "0.1408"
Set-TextValue $ws "E43" "-1.60%"
Set-TextValue $ws "D44" "0.002012"
Set-TextValue $ws "E44" "-6.10%"
Set-TextValue $ws "D45" "0.01131"
Set-TextValue $ws "E45" "17.53%"
Set-TextValue $ws "D46" "0.00006802"
Set-TextValue $ws "E46" "2.10%"
Set-TextValue $ws "D47" "0.00000000754"
Set-TextValue $ws "E47" "0.44%"
Set-TextValue $ws "D48" "0.003044"
Set-TextValue $ws "E48" "3.58%"
Set-TextValue $ws "D49" "0.002280"
Set-TextValue $ws "E49" "-0.98%"
Set-TextValue $ws "D50" "0.00002112"
Set-TextValue $ws "E50" "0.44%"
Set-TextValue $ws "D51" "0.0002011"
Set-TextValue $ws "E51" "0.44%"

Write-Host "Updated 35 rows (70 cells) of Price/Volume(1h) data."
